$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.397.85'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.90'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7138'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.76'
$ws.Range("E6").Value = '  +1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07789'
$ws.Range("E8").Value = '  -3.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3081'
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.99'
$ws.Range("E10").Value = '  +6.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08245'
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7258'
$ws.Range("E12").Value = '  +3.45%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.278'
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.839.96'
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.50'
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.415.07'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.899'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007925'
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.92'
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.30'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.980'
$ws.Range("E22").Value = '  +8.30%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1556'
$ws.Range("E24").Value = '  +8.66%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.61'
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.002'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.31'
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.357'
$ws.Range("E28").Value = '  -3.87%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.480'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.375'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.115'
$ws.Range("E31").Value = '  +2.49%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05268'
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.928'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.199'
$ws.Range("E34").Value = '  +3.64%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7182'
$ws.Range("E35").Value = '  +2.77%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.678'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01861'
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.212.37'
$ws.Range("E38").Value = '  +7.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.709'
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9060'
$ws.Range("E40").Value = '  -2.35%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.063'
$ws.Range("E41").Value = '  +3.83%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.34'
$ws.Range("E42").Value = '  +4.16%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.40'
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5343'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000121'
$ws.Range("E46").Value = '  +3.93%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.761'
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.905'
$ws.Range("E48").Value = '  +9.96%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4317'
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.233'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("B51").Value = 'Frax'
$ws.Range("C51").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9998'
$ws.Range("E51").Value = '  -0.31%  '
